$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Andrew Case's logged time was updated by +0:20 (66h 20m -> 66h 40m).
# Select the cell first so the active selection matches the edit location.
$ws.Range("B4").Select() | Out-Null
$ws.Range("B4").Value = "66h 40m"
